$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1
$ws.Range("C1").Value = "Bucketizer"

# Delete rows 2 and 3 (data rows), leaving only the header row
$ws.Range("A2:M3").EntireRow.Delete()
